# Insert a new data row at row 74 (pushing the existing rows 74..207 down to
# 75..208) and populate it with the new "Jengibre" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).Insert()

$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44868
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100114007
$ws.Range("G74").Value = "Jengibre"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 200
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = 19000
$ws.Range("N74").Value = '$/caja 13 kilos'
$ws.Range("O74").Value = "Perú"
$ws.Range("P74").Value = 1462
$ws.Range("Q74").Value = 13
$ws.Range("R74").Value = "Hortaliza"
